$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1359.54
$ws.Range("F3").Value = 1365.6
$ws.Range("F4").Value = 1372.8
$ws.Range("F5").Value = 1394.4
$ws.Range("F6").Value = 1404.66
